$p = $ppt.ActivePresentation

# --- Slide 12 ("Demo Scenario 2 (network partition)") ---
# Text box "文字方塊 35" (shape 7): rename the candidates from A/B/C to D/E/F
# and drop the now-irrelevant "a1 login, connect to Server 1" line, folding
# its neighbouring result line into the updated D/E tally. The shape has
# auto-fit-height turned on, so resetting .Text recomputes <a:ext cy="..."/>
# to match the new (shorter) paragraph count automatically.
$s12 = $p.Slides.Item(12)
$tb = $s12.Shapes.Item(7)
$lines = @(
    "a1 and a2 login →",
    "Create election : queen(D E F) →",
    "Network partition",
    "-------------------------------------------------------------",
    "a1 vote for D →",
    "a2 vote for E →",
    "Get result respectively →",
    "Network relive →",
    "Get result (D:1, E:1)"
)
$tb.TextFrame.TextRange.Text = [string]::Join("`r", $lines)

# --- Slide 13 ("Demo link") ---
# Content placeholder (shape 2): turn the bare "YouTube : " label into a real
# hyperlink run pointing at the demo video, leaving the label text itself
# untouched.
$s13 = $p.Slides.Item(13)
$link = $s13.Shapes.Item(2)
$tr = $link.TextFrame.TextRange
$priorLen = $tr.Length
$appended = $tr.InsertAfter("https://youtu.be/eROVNXL8Ntw")
$urlRange = $tr.Characters($priorLen + 1, $appended.Length - $priorLen)
$urlRange.ActionSettings.Item(1).Hyperlink.Address = "https://youtu.be/eROVNXL8Ntw"
